# creates.xlsx update — add "The Rolling Stones" (Band ID 781) data block and
# fix the generated SQL text from "insert into creates into(" to
# "insert into creates values(".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the K1 helper text used by the CONCATENATE() formulas
$ws.Range("K1").Value = 'insert into creates values('

# 2. New data block for "The Rolling Stones" (Band ID 781), rows 10-16
$ws.Range("A10").Value = 781
$ws.Range("B10").Value = 3282
$ws.Range("A11").Value = 781
$ws.Range("B11").Value = 3627
$ws.Range("A12").Value = 781
$ws.Range("B12").Value = 7329
$ws.Range("A13").Value = 781
$ws.Range("B13").Value = 7482
$ws.Range("A14").Value = 781
$ws.Range("B14").Value = 5738
$ws.Range("A15").Value = 781
$ws.Range("B15").Value = 6234
$ws.Range("A16").Value = 781
$ws.Range("B16").Value = 7231

# Copy the formatting (center/center/wrap) of the "Reincidentes" label block
# (D6:E8) onto the new label block so we reuse the existing cell style
# instead of minting new ones.
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D10:E16").PasteSpecial(-4122) | Out-Null

# Label for the new block, in D10 (will become merged D10:E16)
$ws.Range("D10").Value = "#The Rolling Stones´s Band ID: 781"

# Merge the label cells for the new block
$ws.Range("D10:E16").Merge()

# SQL-insert helper formulas for the new rows
$ws.Range("I10:I16").Formula = '=CONCATENATE($K$1,A10,$K$2,B10,$L$1)'

# Column K grew wider to fit the new, longer helper text
$ws.Columns.Item(11).ColumnWidth = 24

# Update the active selection to the new block
$ws.Range("A10:A16").Select() | Out-Null
